# Apply the edit described by the diff:
# - Sheet1!A4 value changes from "TBD" to "tba"
# - Selection changes to A4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Name cell for the "Vice President of Pledge Education" row
$ws.Range("A4").Value = "tba"

# Update the selected cell/range to A4, as shown in the diff
$ws.Activate()
$ws.Range("A4").Select()
